$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.968.60'
$ws.Range('E2').Value = '  +0.44%  '

$ws.Range('D3').Value = '1.885.37'
$ws.Range('E3').Value = '  -0.48%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.95'
$ws.Range('E5').Value = '  -2.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.20%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4593'
$ws.Range('E7').Value = '  -2.99%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4069'
$ws.Range('E8').Value = '  +0.72%  '

$ws.Range('E9').Value = '  -0.24%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07968'
$ws.Range('E10').Value = '  -1.58%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9901'
$ws.Range('E11').Value = '  -2.32%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.65'
$ws.Range('E12').Value = '  -2.90%  '

$ws.Range('D13').Value = '1.896.98'
$ws.Range('E13').Value = '  +0.00%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.902'
$ws.Range('E14').Value = '  -2.88%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.047'
$ws.Range('E15').Value = '  -3.88%  '

$ws.Range('E16').Value = '  +0.15%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.32'
$ws.Range('E17').Value = '  -2.91%  '

$ws.Range('E18').Value = '  -2.33%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06546'
$ws.Range('E19').Value = '  -0.84%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.40'
$ws.Range('E20').Value = '  -1.72%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.15%  '

$ws.Range('D22').Value = '29.015.25'
$ws.Range('E22').Value = '  +0.56%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.406'
$ws.Range('E23').Value = '  -2.23%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.22'
$ws.Range('E24').Value = '  +1.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.206'
$ws.Range('E25').Value = '  -2.50%  '

$ws.Range('D26').Value = '2.126.21'
$ws.Range('E26').Value = '  +0.25%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.49'
$ws.Range('E27').Value = '  -2.54%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.58'
$ws.Range('E28').Value = '  -1.70%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.097'
$ws.Range('E29').Value = '  -2.40%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.406'
$ws.Range('E30').Value = '  -2.37%  '

$ws.Range('E31').Value = '  -2.26%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9749'
$ws.Range('E32').Value = '  -2.87%  '

$ws.Range('E33').Value = '  -2.65%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.602'
$ws.Range('E34').Value = '  -1.28%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.401'
$ws.Range('E35').Value = '  +0.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.267'
$ws.Range('E36').Value = '  -2.35%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06039'
$ws.Range('E37').Value = '  -2.46%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02214'
$ws.Range('E38').Value = '  -3.18%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.253'
$ws.Range('E39').Value = '  -3.71%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.177'
$ws.Range('E40').Value = '  -1.11%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9999'
$ws.Range('E41').Value = '  +0.15%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5756'
$ws.Range('E42').Value = '  -4.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1816'
$ws.Range('E43').Value = '  -4.37%  '

$ws.Range('E44').Value = '  -3.10%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.254'
$ws.Range('E45').Value = '  -0.77%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07542'
$ws.Range('E46').Value = '  +4.11%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.241'
$ws.Range('E47').Value = '  +5.59%  '

$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5441'
$ws.Range('E48').Value = '  -2.99%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.91'
$ws.Range('E49').Value = '  -3.06%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.888'
$ws.Range('E50').Value = '  -4.29%  '

$ws.Range('E51').Value = '  -1.57%  '
